$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing "Realization"/"Bring Up"/"Test " rows (15-17) down to
# rows 24-26 by inserting 9 new blank rows at row 15.
$ws.Rows("15:23").Insert()

# Row 15: Source
$ws.Range("B15").Value = "Source"
$ws.Range("D15").Value = 10
$ws.Range("E15").Value = 1
$ws.Range("F15").Formula = "=D15*E15"

# Row 16: Major Compnents BOM
$ws.Range("B16").Value = "Major Compnents BOM"
$ws.Range("F16").Value = 2

# Row 17: HW Block Diagram
$ws.Range("B17").Value = "HW Block Diagram"

# Row 19: Product Architecture (entered before row 18's value, to match
# the original shared-string ordering)
$ws.Range("B19").Value = "Product Architecture"

# Row 20: Calculations
$ws.Range("C20").Value = "Calculations"

# Row 18: Datasheet Research
$ws.Range("B18").Value = "Datasheet Research"

# Row 21: Current Budget
$ws.Range("C21").Value = "Current Budget "

# Row 22: Major Signals List
$ws.Range("C22").Value = "Major Signals List"

# Row 23: Interface Pinouts
$ws.Range("C23").Value = "Interface Pinouts"

# Column width adjustments (autofit to content, as in source)
$ws.Columns("A:C").EntireColumn.AutoFit()
$ws.Columns("I").EntireColumn.AutoFit()

# Restore selection similar to source
$ws.Range("B29").Select() | Out-Null
